$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 346, shifting the existing rows 346-358
# down to 348-360.
$ws.Rows("346:347").Insert()

# Row 346 - new weekly entry (Rosara, 1a nueva(o))
$ws.Range("A346").Value = 3
$ws.Range("B346").Value = "Femacal de La Calera"
$ws.Range("C346").Value = "Coquimbo"
$ws.Range("D346").Value = 44509
$ws.Range("D346").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E346").Value = 5
$ws.Range("F346").Value = 100114001
$ws.Range("G346").Value = "Papa"
$ws.Range("H346").Value = "Rosara"
$ws.Range("I346").Value = "1a nueva(o)"
$ws.Range("J346").Value = 370
$ws.Range("K346").Value = 9000
$ws.Range("L346").Value = 9500
$ws.Range("M346").Value = 9243
$ws.Range("N346").Value = "$/saco 25 kilos"
$ws.Range("O346").Value = "Provincia de Quillota"
$ws.Range("P346").Value = 370
$ws.Range("Q346").Value = 25
$ws.Range("R346").Value = "Hortaliza"

# Row 347 - new weekly entry (Rosara, 2a nueva(o))
$ws.Range("A347").Value = 3
$ws.Range("B347").Value = "Femacal de La Calera"
$ws.Range("C347").Value = "Coquimbo"
$ws.Range("D347").Value = 44509
$ws.Range("D347").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E347").Value = 5
$ws.Range("F347").Value = 100114001
$ws.Range("G347").Value = "Papa"
$ws.Range("H347").Value = "Rosara"
$ws.Range("I347").Value = "2a nueva(o)"
$ws.Range("J347").Value = 180
$ws.Range("K347").Value = 8500
$ws.Range("L347").Value = 8500
$ws.Range("M347").Value = 8500
$ws.Range("N347").Value = "$/saco 25 kilos"
$ws.Range("O347").Value = "Provincia de Quillota"
$ws.Range("P347").Value = 340
$ws.Range("Q347").Value = 25
$ws.Range("R347").Value = "Hortaliza"
